# Update currentAveragePrice / LevePrice / LeveProfit columns (H:N) on several
# leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, reflecting
# refreshed market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 24240608
$ws.Range("I19").Value = 18783576
$ws.Range("J19").Value = 33335664
$ws.Range("K19").Value = 18783576
$ws.Range("L19").Value = 33335664
$ws.Range("M19").Value = -18783401
$ws.Range("N19").Value = -33336014
$ws.Range("H32").Value = 3883.8333
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3883.8333
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3883.8333
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4535.8333
$ws.Range("H33").Value = 355.77777
$ws.Range("I33").Value = 120
$ws.Range("J33").Value = 650.5
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 650.5
$ws.Range("M33").Value = 109
$ws.Range("N33").Value = -1108.5
$ws.Range("H53").Value = 204.72223
$ws.Range("I53").Value = 278.16666
$ws.Range("J53").Value = 168
$ws.Range("K53").Value = 278.16666
$ws.Range("L53").Value = 168
$ws.Range("M53").Value = 358.83334
$ws.Range("N53").Value = -1442
$ws.Range("H76").Value = 3166.5557
$ws.Range("J76").Value = 3139.88
$ws.Range("L76").Value = 3139.88
$ws.Range("N76").Value = -3769.88
$ws.Range("H79").Value = 3166.5557
$ws.Range("J79").Value = 3139.88
$ws.Range("L79").Value = 3139.88
$ws.Range("N79").Value = -5323.88
$ws.Range("H86").Value = 1869.3462
$ws.Range("J86").Value = 1963.6364
$ws.Range("L86").Value = 1963.6364
$ws.Range("N86").Value = -4209.6364
$ws.Range("H89").Value = 1869.3462
$ws.Range("J89").Value = 1963.6364
$ws.Range("L89").Value = 9818.182000000001
$ws.Range("N89").Value = -21050.182
$ws.Range("H116").Value = 3208886.2
$ws.Range("I116").Value = 25643590
$ws.Range("K116").Value = 25643590
$ws.Range("M116").Value = -25640148
$ws.Range("H137").Value = 2666.6538
$ws.Range("I137").Value = 1158.6875
$ws.Range("K137").Value = 3476.0625
$ws.Range("M137").Value = -926.0625
$ws.Range("H138").Value = 3051.5518
$ws.Range("I138").Value = 1020.73914
$ws.Range("J138").Value = 4386.086
$ws.Range("K138").Value = 3062.21742
$ws.Range("L138").Value = 13158.258
$ws.Range("M138").Value = 2077.78258
$ws.Range("N138").Value = -23438.258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1385.52
$ws.Range("I2").Value = 1366.9
$ws.Range("J2").Value = 1460
$ws.Range("K2").Value = 1366.9
$ws.Range("L2").Value = 1460
$ws.Range("M2").Value = -1253.9
$ws.Range("N2").Value = -1686
$ws.Range("H32").Value = 16670870
$ws.Range("I32").Value = 18870678
$ws.Range("J32").Value = 15175.429
$ws.Range("K32").Value = 18870678
$ws.Range("L32").Value = 15175.429
$ws.Range("M32").Value = -18870391
$ws.Range("N32").Value = -15749.429
$ws.Range("H45").Value = 1298.2222
$ws.Range("I45").Value = 947.5
$ws.Range("K45").Value = 947.5
$ws.Range("M45").Value = -570.5
$ws.Range("H63").Value = 2901.8
$ws.Range("J63").Value = 2834.6667
$ws.Range("L63").Value = 2834.6667
$ws.Range("N63").Value = -4206.6667
$ws.Range("H66").Value = 2901.8
$ws.Range("J66").Value = 2834.6667
$ws.Range("L66").Value = 14173.3335
$ws.Range("N66").Value = -21037.3335
$ws.Range("H74").Value = 1322.85
$ws.Range("I74").Value = 1237.1428
$ws.Range("K74").Value = 1237.1428
$ws.Range("M74").Value = -363.1428000000001
$ws.Range("H77").Value = 1322.85
$ws.Range("I77").Value = 1237.1428
$ws.Range("K77").Value = 6185.714
$ws.Range("M77").Value = -1817.714
$ws.Range("H112").Value = 19221.75
$ws.Range("J112").Value = 19221.75
$ws.Range("L112").Value = 19221.75
$ws.Range("N112").Value = -22175.75
$ws.Range("H116").Value = 1385.52
$ws.Range("I116").Value = 1366.9
$ws.Range("J116").Value = 1460
$ws.Range("K116").Value = 1366.9
$ws.Range("L116").Value = 1460
$ws.Range("M116").Value = 927.0999999999999
$ws.Range("N116").Value = -6048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1385.52
$ws.Range("I3").Value = 1366.9
$ws.Range("J3").Value = 1460
$ws.Range("K3").Value = 1366.9
$ws.Range("L3").Value = 1460
$ws.Range("M3").Value = -1252.9
$ws.Range("N3").Value = -1688
$ws.Range("H94").Value = 664.75
$ws.Range("I94").Value = 659.7143
$ws.Range("K94").Value = 659.7143
$ws.Range("M94").Value = -208.7143
$ws.Range("H105").Value = 2498.45
$ws.Range("I105").Value = 1129.1818
$ws.Range("K105").Value = 1129.1818
$ws.Range("M105").Value = 617.8181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2265.1843
$ws.Range("I31").Value = 1709.1072
$ws.Range("J31").Value = 3822.2
$ws.Range("K31").Value = 1709.1072
$ws.Range("L31").Value = 3822.2
$ws.Range("M31").Value = -1414.1072
$ws.Range("N31").Value = -4412.2
$ws.Range("H34").Value = 2265.1843
$ws.Range("I34").Value = 1709.1072
$ws.Range("J34").Value = 3822.2
$ws.Range("K34").Value = 1709.1072
$ws.Range("L34").Value = 3822.2
$ws.Range("M34").Value = -1507.1072
$ws.Range("N34").Value = -4226.2
$ws.Range("H99").Value = 4999.5
$ws.Range("I99").Value = 4999
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 4999
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -3501
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 14997
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12527
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 1538.5676
$ws.Range("I132").Value = 1247.3438
$ws.Range("K132").Value = 3742.0314
$ws.Range("M132").Value = -1212.0314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 613.5641000000001
$ws.Range("I131").Value = 387.04544
$ws.Range("K131").Value = 1161.13632
$ws.Range("M131").Value = 3878.86368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8745.4
$ws.Range("I113").Value = 1581.1428
$ws.Range("J113").Value = 15014.125
$ws.Range("K113").Value = 1581.1428
$ws.Range("L113").Value = 15014.125
$ws.Range("M113").Value = 588.8571999999999
$ws.Range("N113").Value = -19354.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 27000
$ws.Range("J110").Value = 27000
$ws.Range("L110").Value = 27000
$ws.Range("N110").Value = -35180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 351
$ws.Range("I4").Value = 351
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 351
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -238
$ws.Range("N4").ClearContents()
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H122").Value = 3757.7778
$ws.Range("I122").Value = 5333.3335
$ws.Range("J122").Value = 2970
$ws.Range("K122").Value = 16000.0005
$ws.Range("L122").Value = 8910
$ws.Range("M122").Value = -13550.0005
$ws.Range("N122").Value = -13810
$ws.Range("H132").Value = 2334.125
$ws.Range("I132").Value = 1674.0667
$ws.Range("K132").Value = 5022.2001
$ws.Range("M132").Value = -2492.2001
